# Small change to temp contrast drawing (slide 3 of the thesis-drawings deck).
#
# 1) The top-level drawing group on slide 3 was re-grouped by the author
#    (ungroup + regroup), which causes PowerPoint to hand out a fresh
#    shape id/name ("Group 34" -> "Group 7"). We reproduce that id
#    renumbering here: briefly create+delete a throw-away shape to
#    advance the slide's internal id counter to the same point the
#    original authoring session was at, then ungroup/regroup so the
#    resulting group lands on id 8 / "Group 7", exactly like the diff.
# 2) Several child shapes (two connectors and two text boxes) were
#    nudged slightly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$grp = $s.Shapes.Item(1)

$tmp = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tmp.Delete()

$items = $grp.Ungroup()
$grp = $items.Group()
$grp.Name = "Group 7"

$items = $grp.GroupItems

# "Straight Arrow Connector 27" (id 28): flip vertically and resize/move slightly.
$c27 = $items.Item(9)
$c27.VerticalFlip = -1
$c27.Left = 456.7303161621094
$c27.Top = 255.81906127929688
$c27.Width = 145.53582763671875
$c27.Height = 0.2137008011341095

# "Straight Arrow Connector 29" (id 30): move left edge only.
$c29 = $items.Item(11)
$c29.Left = 327.5419006347656

# "TextBox 30" (id 31): move left edge only.
$t30 = $items.Item(12)
$t30.Left = 350.0667724609375

# "TextBox 44" (id 45): move left edge only.
$t44 = $items.Item(15)
$t44.Left = 276.6328430175781
